$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '97.558.36'
Set-TextValue $ws.Range('E2') '  +3.66%  '
Set-TextValue $ws.Range('D3') '3.344.47'
Set-TextValue $ws.Range('E3') '  +8.60%  '
Set-TextValue $ws.Range('E4') '  -0.02%  '
Set-TextValue $ws.Range('D5') '257.55'
Set-TextValue $ws.Range('E5') '  +9.92%  '
Set-TextValue $ws.Range('D6') '620.23'
Set-TextValue $ws.Range('E6') '  +1.86%  '
Set-TextValue $ws.Range('E7') '  +2.32%  '
Set-TextValue $ws.Range('E8') '  +2.09%  '
Set-TextValue $ws.Range('E9') '  +0.05%  '
Set-TextValue $ws.Range('D10') '3.339.34'
Set-TextValue $ws.Range('E10') '  +8.50%  '
Set-TextValue $ws.Range('D11') '0.794'
Set-TextValue $ws.Range('E11') '  -3.34%  '
Set-TextValue $ws.Range('E12') '  +1.78%  '
Set-TextValue $ws.Range('D13') '97.264.26'
Set-TextValue $ws.Range('E13') '  +3.55%  '
Set-TextValue $ws.Range('D14') '35.52'
Set-TextValue $ws.Range('E14') '  +4.52%  '
Set-TextValue $ws.Range('E15') '  +2.58%  '
Set-TextValue $ws.Range('D16') '3.960.77'
Set-TextValue $ws.Range('E16') '  +8.56%  '
Set-TextValue $ws.Range('D17') '5.52'
Set-TextValue $ws.Range('E17') '  +4.72%  '
Set-TextValue $ws.Range('D18') '3.335.87'
Set-TextValue $ws.Range('E18') '  +8.64%  '
Set-TextValue $ws.Range('D19') '3.59'
Set-TextValue $ws.Range('E19') '  -1.22%  '
Set-TextValue $ws.Range('D20') '14.99'
Set-TextValue $ws.Range('E20') '  +3.19%  '
Set-TextValue $ws.Range('D21') '482.04'
Set-TextValue $ws.Range('E21') '  +9.10%  '
Set-TextValue $ws.Range('B22') 'Polkadot'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D22') '5.82'
Set-TextValue $ws.Range('E22') '  +1.47%  '
Set-TextValue $ws.Range('B23') 'PEPE'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D23') '0.0000208'
Set-TextValue $ws.Range('E23') '  +8.58%  '
Set-TextValue $ws.Range('E24') '  +4.73%  '
Set-TextValue $ws.Range('E25') '  +1.83%  '
Set-TextValue $ws.Range('D26') '88.04'
Set-TextValue $ws.Range('E26') '  +3.58%  '
Set-TextValue $ws.Range('D27') '12.10'
Set-TextValue $ws.Range('E27') '  +1.54%  '
Set-TextValue $ws.Range('D28') '3.521.03'
Set-TextValue $ws.Range('E28') '  +8.27%  '
Set-TextValue $ws.Range('E29') '  +0.05%  '
Set-TextValue $ws.Range('D30') '0.184'
Set-TextValue $ws.Range('E30') '  +3.64%  '
Set-TextValue $ws.Range('D31') '0.238'
Set-TextValue $ws.Range('E31') '  -2.81%  '
Set-TextValue $ws.Range('E32') '  -0.65%  '
Set-TextValue $ws.Range('E33') '  +0.63%  '
Set-TextValue $ws.Range('E34') '  +1.17%  '
Set-TextValue $ws.Range('E35') '  +7.03%  '
Set-TextValue $ws.Range('D36') '7.40'
Set-TextValue $ws.Range('E36') '  -4.36%  '
Set-TextValue $ws.Range('E37') '  -4.10%  '
Set-TextValue $ws.Range('D38') '510.48'
Set-TextValue $ws.Range('E38') '  +9.15%  '
Set-TextValue $ws.Range('E39') '  +3.44%  '
Set-TextValue $ws.Range('D40') '24.81'
Set-TextValue $ws.Range('E40') '  +3.46%  '
Set-TextValue $ws.Range('E41') '  +1.20%  '
Set-TextValue $ws.Range('D42') '1.27'
Set-TextValue $ws.Range('E42') '  +0.75%  '
Set-TextValue $ws.Range('D43') '3.31'
Set-TextValue $ws.Range('E43') '  +6.10%  '
Set-TextValue $ws.Range('D44') '3.52'
Set-TextValue $ws.Range('E44') '  -5.47%  '
Set-TextValue $ws.Range('D45') '0.790'
Set-TextValue $ws.Range('E45') '  +16.83%  '
Set-TextValue $ws.Range('E46') '  +0.01%  '
Set-TextValue $ws.Range('D47') '160.94'
Set-TextValue $ws.Range('E47') '  +0.73%  '
Set-TextValue $ws.Range('E48') '  +3.35%  '
Set-TextValue $ws.Range('E49') '  +6.69%  '
Set-TextValue $ws.Range('D50') '45.52'
Set-TextValue $ws.Range('E50') '  +4.10%  '
Set-TextValue $ws.Range('D51') '4.52'
Set-TextValue $ws.Range('E51') '  +5.29%  '
